# Apply updated crypto price/volume figures (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.514.68"
$ws.Range("E2").Value = "  +2.51%  "
$ws.Range("D3").Value = "2.193.22"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "'258.13"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.39%  "
$ws.Range("D6").Value = "'83.31"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +11.47%  "
$ws.Range("D7").Value = "'0.618"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.09%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").Value = "'0.599"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +3.17%  "
$ws.Range("D10").Value = "'44.81"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +11.59%  "
$ws.Range("E11").Value = "  +0.69%  "
$ws.Range("D12").Value = "'7.11"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +5.24%  "
$ws.Range("E13").Value = "  +2.54%  "
$ws.Range("D14").Value = "2.523.56"
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("D15").Value = "'14.32"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.01%  "
$ws.Range("D16").Value = "2.181.78"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").Value = "'0.782"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.77%  "
$ws.Range("D18").Value = "43.432.48"
$ws.Range("E18").Value = "  +2.52%  "
$ws.Range("E19").Value = "  +1.50%  "
$ws.Range("D20").Value = "'69.58"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.76%  "
$ws.Range("D21").Value = "'5.89"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.80%  "
$ws.Range("D22").Value = "'2.36"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +12.20%  "
$ws.Range("D23").Value = "'231.09"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.99%  "
$ws.Range("D24").Value = "'8.93"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -6.03%  "
$ws.Range("D25").Value = "'0.999"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("D26").Value = "'10.60"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.76%  "
$ws.Range("E27").Value = "  +0.69%  "
$ws.Range("D28").Value = "'39.53"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +7.04%  "
$ws.Range("E29").Value = "  +2.91%  "
$ws.Range("E30").Value = "  +3.11%  "
$ws.Range("D31").Value = "'174.12"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.28%  "
$ws.Range("D32").Value = "'20.35"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.74%  "
$ws.Range("D33").Value = "'0.0859"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +5.29%  "
$ws.Range("D34").Value = "'5.32"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.55%  "
$ws.Range("E35").Value = "  +2.32%  "
$ws.Range("E36").Value = "  +4.88%  "
$ws.Range("E37").Value = "  +7.81%  "
$ws.Range("E38").Value = "  +7.56%  "
$ws.Range("D39").Value = "'12.45"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +3.45%  "
$ws.Range("E40").Value = "  +9.19%  "
$ws.Range("E41").Value = "  +1.67%  "
$ws.Range("D42").Value = "'62.94"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +7.48%  "
$ws.Range("D43").Value = "'5.46"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +5.76%  "
$ws.Range("D44").Value = "'0.198"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.83%  "
$ws.Range("D46").Value = "'0.0973"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.38%  "
$ws.Range("D47").Value = "'99.40"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.02%  "
$ws.Range("E48").Value = "  +5.29%  "
$ws.Range("E49").Value = "  +1.89%  "
$ws.Range("D50").Value = "'0.437"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -4.12%  "
$ws.Range("E51").Value = "  +11.95%  "
